# Auto-generated edit script applying cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.683.33"
$ws.Range("E2").Value = "  -0.37%  "
$ws.Range("D3").Value = "2.309.44"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "538.45"
$ws.Range("E5").Value = "  -2.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.38"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  +2.24%  "
$ws.Range("D9").Value = "2.308.36"
$ws.Range("E9").Value = "  +0.17%  "
$ws.Range("E10").Value = "  -1.41%  "
$ws.Range("E11").Value = "  -0.73%  "
$ws.Range("E13").Value = "  -0.81%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.83"
$ws.Range("E14").Value = "  -1.06%  "
$ws.Range("D15").Value = "2.719.91"
$ws.Range("E15").Value = "  +0.45%  "
$ws.Range("D16").Value = "58.579.84"
$ws.Range("E16").Value = "  -0.46%  "
$ws.Range("E17").Value = "  -0.54%  "
$ws.Range("D18").Value = "2.307.54"
$ws.Range("E18").Value = "  +2.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.60"
$ws.Range("E19").Value = "  -0.99%  "
$ws.Range("E20").Value = "  -3.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "315.94"
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("E22").Value = "  +1.53%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.19"
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.171"
$ws.Range("E25").Value = "  -0.53%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.95"
$ws.Range("E27").Value = "  -2.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.29"
$ws.Range("E28").Value = "  -1.10%  "
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.73"
$ws.Range("E29").Value = "  -1.94%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "171.27"
$ws.Range("E30").Value = "  +0.84%  "
$ws.Range("D31").Value = "0.0₃0731"
$ws.Range("E31").Value = "  -0.27%  "
$ws.Range("E32").Value = "  -0.30%  "
$ws.Range("E33").Value = "  +0.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.387"
$ws.Range("E34").Value = "  +0.67%  "
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.93"
$ws.Range("E35").Value = "  +0.67%  "
$ws.Range("B36").Value = "USDe"
$ws.Range("C36").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("E38").Value = "  -0.34%  "
$ws.Range("E39").Value = "  +0.35%  "
$ws.Range("E40").Value = "  -0.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "291.08"
$ws.Range("E41").Value = "  -4.97%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "141.28"
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.47"
$ws.Range("E43").Value = "  +0.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0954"
$ws.Range("E44").Value = "  -0.15%  "
$ws.Range("E45").Value = "  -1.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.558"
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "18.37"
$ws.Range("E47").Value = "  -2.58%  "
$ws.Range("E48").Value = "  -2.63%  "
$ws.Range("E49").Value = "  -0.62%  "
$ws.Range("E50").Value = "  +0.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.52"
$ws.Range("E51").Value = "  +0.02%  "
